$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>
#    markers that bracket the run "3)Design" (paragraph 2). There is no direct
#    object-model call for proofing marks, so we rebuild that single
#    paragraph's contents via InsertXML from a hand-built OOXML fragment that
#    preserves the paragraph/run formatting and rsid attributes but omits the
#    proofErr elements.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$designParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document ' + $wNs + '><w:body>' +
    '<w:p w:rsidR="00F23865" w:rsidRPr="0003090B" w:rsidRDefault="0003090B">' +
    '<w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="0003090B"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>3)Design</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$designPara = $d.Paragraphs.Item(2)
$designPara.Range.InsertXML($designParaXml)

# ---------------------------------------------------------------------------
# 2. Remove the _GoBack bookmark that sits in the paragraph holding Picture 14,
#    remembering which paragraph it was in so we can find our way to the
#    5th blank paragraph after it (step 3) regardless of exact numbering.
# ---------------------------------------------------------------------------
$picture14ParaIndex = 8
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBackRange = $goBack.Range
    $goBackPara = $goBackRange.Paragraphs.Item(1)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $goBackPara.Range.Start) {
            $picture14ParaIndex = $i
        }
    }
    $goBack.Delete()
}

# ---------------------------------------------------------------------------
# 3. Add the text "Git is Imp " into the 5th blank paragraph following the
#    Picture 14 paragraph (previously completely empty).
# ---------------------------------------------------------------------------
$gitParaIndex = $picture14ParaIndex + 5
$gitPara = $d.Paragraphs.Item($gitParaIndex)
$gitPara.Range.InsertBefore("Git is Imp ")
